# Negative_Manifest.xlsx update: append rows 26-39 (n25..n38), keep one
# trailing blank row (40), update dimension/selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# subject_id, file_name, training_subject, feedback_1_id, classification
$rows = @(
    @("n25","n25_IMG_3070 - Copy.jpeg","True","no_meltpatch","negative"),
    @("n26","n26_IMG_3070.jpeg","True","no_meltpatch","negative"),
    @("n27","n27_IMG_3070HorFlip - Copy.jpeg","True","no_meltpatch","negative"),
    @("n28","n28_IMG_3070HorFlip.jpeg","True","no_meltpatch","negative"),
    @("n29","n29_IMG_3070HorVertFlip.jpeg","True","no_meltpatch","negative"),
    @("n30","n30_IMG_3070VertFlip.jpeg","True","no_meltpatch","negative"),
    @("n31","n31_IMG_3072.jpeg","True","no_meltpatch","negative"),
    @("n32","n32_IMG_3072HorFlip.jpeg","True","no_meltpatch","negative"),
    @("n33","n33_IMG_3072HorVertFlip.jpeg","True","no_meltpatch","negative"),
    @("n34","n34_IMG_3072VertFlip.jpeg","True","no_meltpatch","negative"),
    @("n35","n35_IMG_3073.jpeg","True","no_meltpatch","negative"),
    @("n36","n36_IMG_3073HorFlip.jpeg","True","no_meltpatch","negative"),
    @("n37","n37_IMG_3073HorVertFlip.jpeg","True","no_meltpatch","negative"),
    @("n38","n38_IMG_3073VertFlip.jpeg","True","no_meltpatch","negative")
)

$startRow = 26
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]

    # Force the literal text "True" (not a boolean) via the classic
    # apostrophe text-prefix, then reset the style so no quotePrefix
    # formatting artifact is left behind on the cell.
    $ws.Cells.Item($r, 3).Value = "'" + $vals[2]
    $ws.Cells.Item($r, 3).Style = "Normal"

    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
}

$lastRow = $startRow + $rows.Count - 1

# Update selection to match the authored state: activeCell at the first
# new row, selection spanning the newly written block.
[void]$ws.Range("A" + $startRow + ":E" + $lastRow).Select()

Write-Output "done"
